$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new header row at the top of the data (row 1) with column labels
$headers = @("img1x", "img1y", "img2x", "img2y", "img3x", "img3y", "img4x", "img4y")

for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = 2 + $i  # B = 2
    $ws.Cells.Item(1, $col).Value = $headers[$i]
}

# Update selection to reflect the state captured in the saved workbook
$ws.Range("H15").Select()
